$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.320.37'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.666.82'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.43'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5302'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.009'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2643'
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06359'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.92'
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07848'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.529'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '1.677.31'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').Value = '1.895.93'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5602'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '0.0₅8136'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.77'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = '26.336.36'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.717'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '198.55'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.27'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.055'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.20'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1212'
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.227'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.513'
$ws.Range('E29').Value = '  +2.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05881'
$ws.Range('E30').Value = '  +1.48%  '
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.535'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.318'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.833'
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9603'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5799'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01616'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.959'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').Value = '1.074.13'
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8577'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.008'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.81'
$ws.Range('E44').Value = '  -1.84%  '
$ws.Range('D45').Value = '1.806.72'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.48'
$ws.Range('E46').Value = '  +2.47%  '
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4412'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₈104'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05155'
$ws.Range('E51').Value = '  -0.07%  '
